$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024701
$ws.Range("N2").Value = 9.074103000000001
$ws.Range("O2").Value = 0.1596375877334842
$ws.Range("P2").Value = 0.1596375877334843
$ws.Range("Q2").Value = 1.444637526946667
$ws.Range("R2").Value = 13.00173774252
$ws.Range("S2").Value = 0.1596375877334842
$ws.Range("T2").Value = 0.1596375877334843

# Row 3
$ws.Range("O3").Value = 0.6072559333217162
$ws.Range("P3").Value = 0.6072559333217163
$ws.Range("S3").Value = 0.6072559333217162
$ws.Range("T3").Value = 0.6072559333217163

# Row 4
$ws.Range("M4").Value = 4.368554666666666
$ws.Range("N4").Value = 13.105664
$ws.Range("O4").Value = 0.2305634602787257
$ws.Range("P4").Value = 0.2305634602787257
$ws.Range("Q4").Value = 2.086479956195556
$ws.Range("R4").Value = 18.77831960576
$ws.Range("S4").Value = 0.2305634602787257
$ws.Range("T4").Value = 0.2305634602787257

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04818333333333333
$ws.Range("N5").Value = 0.14455
$ws.Range("O5").Value = 0.002543018666073676
$ws.Range("P5").Value = 0.002543018666073677
$ws.Range("Q5").Value = 0.02301300244444445
$ws.Range("R5").Value = 0.207117022
$ws.Range("S5").Value = 0.002543018666073676
$ws.Range("T5").Value = 0.002543018666073677
